$d = $word.ActiveDocument
$vt = [char]11

# 1. Title
$r1 = $d.Content.Find.Execute("Interstellar Explorations: Delving into the Cosmic Abyss", $true, $false, $false, $false, $false, $true, 1, $false, "Time Travel: A Glimpse into Past, Present, and Future", 2)
Write-Output "title: $r1"

# 2. Author name
$r2 = $d.Content.Find.Execute("Dr. Luna A. Hendrix", $true, $false, $false, $false, $false, $true, 1, $false, "Sandra Mandell", 2)
Write-Output "name: $r2"

# 3. Email local-part
$r3 = $d.Content.Find.Execute("luna", $true, $false, $false, $false, $false, $true, 1, $false, "Sandra_mandel@unitedassociatedschools", 2)
Write-Output "email1: $r3"

# 4. Email domain
$r4 = $d.Content.Find.Execute("hendrix21@universeconnect.space", $true, $false, $false, $false, $false, $true, 1, $false, "edu", 2)
Write-Output "email2: $r4"

# 5. Big body paragraph
$old5 = "In the vast expanse of the cosmos, humanity's innate curiosity has propelled us to explore the depths of the unknown. Our desire to uncover the secrets of distant worlds, unravel the mysteries of extraterrestrial life, and grasp the grandeur of the universe knows no bounds. Interstellar explorations have emerged as a captivating realm of scientific endeavor, driven by the indomitable spirit of discovery. Embarking on these epic journeys, we transcend the limits of our earthly existence and embrace the limitless possibilities that lie beyond our planet." + $vt + $vt + "As we venture into the cosmic abyss, we encounter awe-inspiring celestial bodies, each possessing unique characteristics and captivating phenomena. From the fiery hearts of stars and the enigmatic rings of gas giants to the vibrant dance of nebulas and the perplexing properties of black holes, the universe unveils its boundless wonders. Through meticulous observations, cutting-edge technologies, and ingenious space missions, we unravel the mysteries that have long tantalized our imaginations." + $vt + $vt + "Interstellar explorations challenge our understanding of physics, relativity, and the very nature of space and time. As we delve deeper into the cosmos, we encounter gravitational forces, cosmic radiation, and phenomena that defy conventional wisdom. These challenges fuel scientific breakthroughs, pushing the boundaries of human knowledge and propelling us towards a profound comprehension of the fundamental laws that govern the universe."
$new5 = "Imagine traveling through time; visiting ancient civilizations, witnessing pivotal historical events, or experiencing the marvels of the future. While time travel remains a captivating fantasy, it prompts thought-provoking discussions on the nature of time, causality, and the human quest for exploration. In this essay, we delve into the concept of time travel, exploring its possibilities within the realms of theory, science fiction, and pop culture." + $vt + $vt + "Visualize a world where you could journey back to the Renaissance to witness the birth of great art and literature, experience the thrill of a dinosaur-populated prehistoric era, or delve into tomorrow's scientific advancements. Time travel captivates imaginations and ignites conversations about the fluidity of time, the interconnectedness of past, present, and future, and the pursuit of knowledge. However, the concept of time travel is bound by scientific limitations, posing perplexing questions about paradoxes, causality, and the intricate tapestry of temporal threads." + $vt + $vt + "Exploring the subject further, we find a multitude of works in science fiction and literature that delve into the complexities of time travel. Authors like HG. Wells and Ray Bradbury construct captivating worlds where time machines enable voyages across ages, exploring the consequences and implications of altering history. These fictional accounts mirror the human fascination with the unknown and the irresistible allure of unraveling the mysteries of time. They offer a glimpse into alternate realities, inspiring us to ponder the infinite possibilities that exist beyond the constructs of our perceived linear existence." + $vt + $vt + "Body:" + $vt + $vt + "Time travel remains a debated concept within the scientific community. Theories like Einstein's Theory of Relativity suggest that time is not absolute, and that bending spacetime can potentially allow for movement through time. However, the practical implementation of such theories is still beyond our current scientific capabilities, leaving time travel firmly entrenched in the realms of speculation. Despite the scientific challenges, the idea of time travel sparks creativity and drives technological advancements." + $vt + $vt + "The concept of time travel also provokes numerous philosophical conundrums. The grandfather paradox, where an individual travels back in time and prevents their grandparent's meeting, poses questions about causality and the stability of historical events. Questions arise: Can the past be changed? Does changing the past ripple through time, altering the present and future? These paradoxes challenge our understanding of cause and effect, pushing us to contemplate the intricate relationship between different points in time." + $vt + $vt + "Furthermore, time travel ignites conversations about knowledge, control, and the potential consequences of altering the course of history. If we could travel back in time, would we correct past mistakes? How would changing one event affect subsequent historical developments? Would it lead to a utopian future or a disastrous dystopia? These contemplations delve into profound questions about human agency, responsibility, and our role in shaping the world around us."
$r5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
Write-Output "body: $r5"

# 6. Summary body paragraph
$old7 = "Interstellar explorations epitomize the audacity of human ambition, driven by our insatiable quest for knowledge and the allure of the unknown. By embarking on these extraordinary journeys, we unravel the mysteries of distant worlds, encounter celestial spectacles, and challenge our understanding of the fundamental laws of the universe. Each mission expands our cosmic horizons, unveiling the vastness and complexity of the cosmos while inspiring future generations to dream big and strive for the impossible."
$new7 = "The concept of time travel captivates imaginations and prompts profound conversations about the nature of time, causality, and the allure of the unknown. While its scientific feasibility remains a topic of debate, time travel continues to inspire creative endeavors, such as science fiction literature and film. Through theoretical musings, fictional narratives, and philosophical inquiries, time travel transcends its speculative nature and provides a lens through which we can examine our own existence and the complexities of the universe. It is a testament to the boundless nature of human imagination and our enduring quest for knowledge and comprehension of the mysteries that lie beyond the limits of our current understanding."
$r7 = $d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2)
Write-Output "summary: $r7"

# 7. Trailing empty paragraph
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertParagraphAfter()
Write-Output "paragraphs: $($d.Paragraphs.Count)"

Write-Output $d.Content.Text
